# Apply edits to Sheet1 per the target diff (ticker list refresh + 7 new rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column B (ticker, col 2) for rows 2..37 ----
$bValues = @{
    2  = "NSE:ALANKIT"
    3  = "NSE:APOLLOHOSP"
    4  = "NSE:BAJAJHCARE"
    5  = "NSE:CAMPUS"
    6  = "NSE:CAPLIPOINT"
    7  = "NSE:CEATLTD"
    8  = "NSE:CENTENKA"
    9  = "NSE:CHAMBLFERT"
    10 = "NSE:CONSOFINVT"
    11 = "NSE:DHARMAJ"
    12 = "NSE:DRCSYSTEMS"
    13 = "NSE:FACT"
    14 = "NSE:FIBERWEB"
    15 = "NSE:FMGOETZE"
    16 = "NSE:GEECEE"
    17 = "NSE:GENUSPOWER"
    18 = "NSE:HARSHA"
    19 = "NSE:IFCI"
    20 = "NSE:INDOCO"
    21 = "NSE:ITDC"
    22 = "NSE:ITI"
    23 = "NSE:JMFINANCIL"
    24 = "NSE:KIRIINDUS"
    25 = "NSE:KSL"
    26 = "NSE:MADRASFERT"
    27 = "NSE:MANORAMA"
    28 = "NSE:MOL"
    29 = "NSE:NAGAFERT"
    30 = "NSE:NAHARINDUS"
    31 = "NSE:NAVINFLUOR"
    32 = "NSE:NFL"
    33 = "NSE:PANSARI"
    34 = "NSE:PATINTLOG"
    35 = "NSE:POLYPLEX"
    36 = "NSE:RCF"
    37 = "NSE:ROSSARI"
}

foreach ($r in $bValues.Keys) {
    $ws.Cells.Item($r, 2).Value = $bValues[$r]
}

# ---- Column C (col 3) updates for rows 2..6 ----
$cValues = @{
    2 = "NSE:ICICIPRULI"
    3 = "NSE:KALAMANDIR"
    4 = "NSE:LUPIN"
    5 = "NSE:METROBRAND"
    6 = "NSE:RAYMOND"
}
foreach ($r in $cValues.Keys) {
    $ws.Cells.Item($r, 3).Value = $cValues[$r]
}

# ---- Column D (col 4): clear D2 and D3 (previously had tickers) ----
$ws.Cells.Item(2, 4).ClearContents()
$ws.Cells.Item(3, 4).ClearContents()

# ---- Column E (col 5): E2 updated ----
$ws.Cells.Item(2, 5).Value = "NSE:ABFRL"

# ---- Column F (col 6): F2 and F3 updated ----
$ws.Cells.Item(2, 6).Value = "NSE:APOLLOHOSP"
$ws.Cells.Item(3, 6).Value = "NSE:NAVINFLUOR"

# ---- New rows 31..37: column A index numbers, styled like the existing A column ----
$aValues = @{
    31 = 29
    32 = 30
    33 = 31
    34 = 32
    35 = 33
    36 = 34
    37 = 35
}

# Reuse the formatting of the existing A column (bold, bordered, centered)
# by copying the format from row 30 down into the new rows, then write values.
$fmtSource = $ws.Cells.Item(30, 1)
$fmtSource.Copy()
foreach ($r in $aValues.Keys) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = $false

foreach ($r in $aValues.Keys) {
    $ws.Cells.Item($r, 1).Value = $aValues[$r]
}

# Update the sheet dimension to reflect the new extent.
$ws.UsedRange | Out-Null
